$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update the "Current page (bookmark)" for 'Researching Information Systems and Computing'
# (row 11) from 129 to 132, reflecting continued reading progress.
$ws.Range("C11").Value = 132

# Move the active selection/cell as recorded in the saved view (C15 -> C21).
$ws.Range("C21").Select()
